$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 107.5
$ws.Range("I4").Value = 116.666664
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 116.666664
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = -2.666663999999997
$ws.Range("N4").Value = -308

$ws.Range("H33").Value = 114.75
$ws.Range("I33").Value = 114.75
$ws.Range("K33").Value = 114.75
$ws.Range("M33").Value = 114.25

$ws.Range("H98").Value = 230051
$ws.Range("I98").Value = 255992.19
$ws.Range("J98").Value = 1768.6
$ws.Range("K98").Value = 255992.19
$ws.Range("L98").Value = 1768.6
$ws.Range("M98").Value = -254494.19
$ws.Range("N98").Value = -4764.6

$ws.Range("H100").Value = 166666670
$ws.Range("I100").Value = 166666670
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 166666670
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -166666129
$ws.Range("N100").ClearContents()

$ws.Range("H112").Value = 7576985.5
$ws.Range("I112").Value = 200
$ws.Range("J112").Value = 8022678.5
$ws.Range("K112").Value = 600
$ws.Range("L112").Value = 24068035.5
$ws.Range("M112").Value = 508
$ws.Range("N112").Value = -24070251.5

$ws.Range("H122").Value = 230051
$ws.Range("I122").Value = 255992.19
$ws.Range("J122").Value = 1768.6
$ws.Range("K122").Value = 767976.5700000001
$ws.Range("L122").Value = 5305.799999999999
$ws.Range("M122").Value = -765526.5700000001
$ws.Range("N122").Value = -10205.8

$ws.Range("H137").Value = 22223568
$ws.Range("I137").Value = 27028020
$ws.Range("J137").Value = 2984.875
$ws.Range("K137").Value = 81084060
$ws.Range("L137").Value = 8954.625
$ws.Range("M137").Value = -81081510
$ws.Range("N137").Value = -14054.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1043.9412
$ws.Range("I2").Value = 723.26086
$ws.Range("J2").Value = 1714.4546
$ws.Range("K2").Value = 723.26086
$ws.Range("L2").Value = 1714.4546
$ws.Range("M2").Value = -610.26086
$ws.Range("N2").Value = -1940.4546

$ws.Range("H110").Value = 493.33334
$ws.Range("I110").Value = 390
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 390
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 1655
$ws.Range("N110").Value = -4790

$ws.Range("H116").Value = 1043.9412
$ws.Range("I116").Value = 723.26086
$ws.Range("J116").Value = 1714.4546
$ws.Range("K116").Value = 723.26086
$ws.Range("L116").Value = 1714.4546
$ws.Range("M116").Value = 1570.73914
$ws.Range("N116").Value = -6302.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1043.9412
$ws.Range("I3").Value = 723.26086
$ws.Range("J3").Value = 1714.4546
$ws.Range("K3").Value = 723.26086
$ws.Range("L3").Value = 1714.4546
$ws.Range("M3").Value = -609.26086
$ws.Range("N3").Value = -1942.4546

$ws.Range("H86").Value = 10592.363
$ws.Range("I86").Value = 1667
$ws.Range("J86").Value = 21302.8
$ws.Range("K86").Value = 1667
$ws.Range("L86").Value = 21302.8
$ws.Range("M86").Value = -544
$ws.Range("N86").Value = -23548.8

$ws.Range("H89").Value = 10592.363
$ws.Range("I89").Value = 1667
$ws.Range("J89").Value = 21302.8
$ws.Range("K89").Value = 8335
$ws.Range("L89").Value = 106514
$ws.Range("M89").Value = -2719
$ws.Range("N89").Value = -117746

$ws.Range("H94").Value = 747.9048
$ws.Range("I94").Value = 780.3
$ws.Range("J94").Value = 100
$ws.Range("K94").Value = 780.3
$ws.Range("L94").Value = 100
$ws.Range("M94").Value = -329.3
$ws.Range("N94").Value = -1002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 31251770
$ws.Range("I86").Value = 55557188
$ws.Range("J86").Value = 1942.5714
$ws.Range("K86").Value = 55557188
$ws.Range("L86").Value = 1942.5714
$ws.Range("M86").Value = -55556065
$ws.Range("N86").Value = -4188.5714

$ws.Range("H89").Value = 31251770
$ws.Range("I89").Value = 55557188
$ws.Range("J89").Value = 1942.5714
$ws.Range("K89").Value = 277785940
$ws.Range("L89").Value = 9712.857
$ws.Range("M89").Value = -277780324
$ws.Range("N89").Value = -20944.857

$ws.Range("H135").Value = 41038.46
$ws.Range("J135").Value = 41038.46
$ws.Range("L135").Value = 41038.46
$ws.Range("N135").Value = -51178.46

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 806.1667
$ws.Range("I107").Value = 1172.4445
$ws.Range("J107").Value = 439.8889
$ws.Range("K107").Value = 1172.4445
$ws.Range("L107").Value = 439.8889
$ws.Range("M107").Value = 747.5554999999999
$ws.Range("N107").Value = -4279.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4784.2856
$ws.Range("I40").Value = 3490
$ws.Range("J40").Value = 4883.846
$ws.Range("K40").Value = 3490
$ws.Range("L40").Value = 4883.846
$ws.Range("M40").Value = -3354
$ws.Range("N40").Value = -5155.846

$ws.Range("H55").Value = 376.26315
$ws.Range("I55").Value = 305.0909
$ws.Range("J55").Value = 474.125
$ws.Range("K55").Value = 305.0909
$ws.Range("L55").Value = 474.125
$ws.Range("M55").Value = -132.0909
$ws.Range("N55").Value = -820.125

$ws.Range("H100").Value = 2941.1765
$ws.Range("I100").Value = 2400
$ws.Range("J100").Value = 3057.1428
$ws.Range("K100").Value = 2400
$ws.Range("L100").Value = 3057.1428
$ws.Range("M100").Value = -1859
$ws.Range("N100").Value = -4139.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 3333.3333
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2351
$ws.Range("N40").Value = -5298

$ws.Range("H81").Value = 503188.78
$ws.Range("I81").Value = 1334383.4
$ws.Range("J81").Value = 4472
$ws.Range("K81").Value = 2668766.8
$ws.Range("L81").Value = 8944
$ws.Range("M81").Value = -2667705.8
$ws.Range("N81").Value = -11066

$ws.Range("H84").Value = 503188.78
$ws.Range("I84").Value = 1334383.4
$ws.Range("J84").Value = 4472
$ws.Range("K84").Value = 13343834
$ws.Range("L84").Value = 44720
$ws.Range("M84").Value = -13338530
$ws.Range("N84").Value = -55328

$ws.Range("H107").Value = 3268578.5
$ws.Range("I107").Value = 3268578.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 9805735.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -9803815.5
$ws.Range("N107").ClearContents()

$ws.Range("H123").Value = 26034.096
$ws.Range("J123").Value = 26034.096
$ws.Range("L123").Value = 26034.096
$ws.Range("N123").Value = -35834.09600000001
